$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "42.708.47"
$ws.Range("E2").Value = "  -0.67%  "
$ws.Range("D3").Value = "2.531.86"
$ws.Range("E3").Value = "  -0.86%  "
$ws.Range("E4").Value = "  -0.05%  "
$ws.Range("D5").Value = "'316.28"
$ws.Range("E5").Value = "  +3.68%  "
$ws.Range("D6").Value = "'95.65"
$ws.Range("E6").Value = "  -2.69%  "
$ws.Range("E7").Value = "  +0.46%  "
$ws.Range("E8").Value = "  +0.04%  "
$ws.Range("E9").Value = "  -1.51%  "
$ws.Range("D10").Value = "'36.54"
$ws.Range("E10").Value = "  -1.42%  "
$ws.Range("D11").Value = "'0.0813"
$ws.Range("E11").Value = "  -1.81%  "
$ws.Range("D12").Value = "'7.76"
$ws.Range("E12").Value = "  -0.24%  "
$ws.Range("E13").Value = "  -1.29%  "
$ws.Range("D14").Value = "2.920.32"
$ws.Range("E14").Value = "  -0.89%  "
$ws.Range("D15").Value = "'15.68"
$ws.Range("E15").Value = "  +3.44%  "
$ws.Range("D16").Value = "2.562.27"
$ws.Range("E16").Value = "  -2.48%  "
$ws.Range("E17").Value = "  -1.74%  "
$ws.Range("D18").Value = "42.749.56"
$ws.Range("E18").Value = "  -0.65%  "
$ws.Range("D19").Value = "'13.13"
$ws.Range("E19").Value = "  -4.90%  "
$ws.Range("D20").Value = "'6.66"
$ws.Range("E20").Value = "  +0.95%  "
$ws.Range("D21").Value = "0.0₃0972"
$ws.Range("E21").Value = "  -2.10%  "
$ws.Range("D22").Value = "'71.36"
$ws.Range("E22").Value = "  -0.87%  "
$ws.Range("D23").Value = "'254.55"
$ws.Range("E23").Value = "  -0.29%  "
$ws.Range("E24").Value = "  +0.76%  "
$ws.Range("E25").Value = "  -1.92%  "
$ws.Range("D26").Value = "'27.61"
$ws.Range("E26").Value = "  -1.78%  "
$ws.Range("E27").Value = "  -0.50%  "
$ws.Range("E28").Value = "  +12.54%  "
$ws.Range("D29").Value = "'39.75"
$ws.Range("E29").Value = "  +4.87%  "
$ws.Range("D30").Value = "'10.09"
$ws.Range("E30").Value = "  -1.75%  "
$ws.Range("D31").Value = "'5.93"
$ws.Range("E31").Value = "  -4.37%  "
$ws.Range("D32").Value = "'156.01"
$ws.Range("E32").Value = "  -1.51%  "
$ws.Range("D33").Value = "'20.07"
$ws.Range("E33").Value = "  +2.19%  "
$ws.Range("D34").Value = "'2.14"
$ws.Range("E34").Value = "  -1.00%  "
$ws.Range("D35").Value = "'3.36"
$ws.Range("E35").Value = "  +0.93%  "
$ws.Range("D36").Value = "'0.0791"
$ws.Range("E36").Value = "  -1.69%  "
$ws.Range("E37").Value = "  -0.76%  "
$ws.Range("B38").Value = "EnergySwap"
$ws.Range("C38").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D38").Value = "'25.31"
$ws.Range("E38").Value = "  -0.99%  "
$ws.Range("B39").Value = "Kaspa"
$ws.Range("C39").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D39").Value = "'0.113"
$ws.Range("E39").Value = "  -3.21%  "
$ws.Range("E40").Value = "  -0.10%  "
$ws.Range("E41").Value = "  +7.92%  "
$ws.Range("D42").Value = "'3.40"
$ws.Range("E42").Value = "  -1.85%  "
$ws.Range("E43").Value = "  -1.71%  "
$ws.Range("D44").Value = "'0.0303"
$ws.Range("E44").Value = "  -1.21%  "
$ws.Range("E45").Value = "  -0.04%  "
$ws.Range("D46").Value = "2.044.33"
$ws.Range("E46").Value = "  -2.35%  "
$ws.Range("D47").Value = "'85.84"
$ws.Range("E47").Value = "  -1.24%  "
$ws.Range("D48").Value = "'8.89"
$ws.Range("E48").Value = "  -1.38%  "
$ws.Range("D49").Value = "'74.89"
$ws.Range("E49").Value = "  -0.65%  "
$ws.Range("D50").Value = "2.775.69"
$ws.Range("E50").Value = "  -0.98%  "
$ws.Range("E51").Value = "  -1.08%  "
